# Append the next 3 days of GSC video-indexing export data to the
# "Chart" sheet (sheet1), mirroring the existing row pattern:
#   A: date (text, "yyyy-MM-dd")
#   B: "No video indexed" count
#   C: "Video indexed" count
#   D: Impressions

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRows = @(
    @("2025-12-09", 23, 1, 0),
    @("2025-12-10", 23, 1, 0),
    @("2025-12-11", 23, 1, 0)
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

foreach ($row in $newRows) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = "'" + $row[0]
    $ws.Cells.Item($lastRow, 2).Value = $row[1]
    $ws.Cells.Item($lastRow, 3).Value = $row[2]
    $ws.Cells.Item($lastRow, 4).Value = $row[3]
}
